# Updated cryptos list: refreshed Price (D) and Volume(1h) (E) columns
# with the latest scrape values; row 51 swapped from SynthetixNetwork (SNX)
# to RenderToken (RNDR) with its new link, price and change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text, preserving the cell's existing style,
# so price strings like '214.53' or '0.500' are not reinterpreted as
# numbers (which would drop meaningful trailing zeros / formatting).
function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue 'D2' '25.751.22'
$ws.Range('E2').Value = '  -0.22%  '
Set-TextValue 'D3' '1.630.00'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '214.53'
$ws.Range('E5').Value = '  -0.35%  '
Set-TextValue 'D6' '0.500'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.86%  '
Set-TextValue 'D9' '0.0631'
$ws.Range('E9').Value = '  -1.59%  '
Set-TextValue 'D10' '19.44'
$ws.Range('E10').Value = '  -2.28%  '
Set-TextValue 'D11' '0.0794'
$ws.Range('E11').Value = '  +0.98%  '
Set-TextValue 'D12' '4.25'
$ws.Range('E12').Value = '  +0.15%  '
Set-TextValue 'D13' '1.856.54'
$ws.Range('E13').Value = '  -0.25%  '
Set-TextValue 'D14' '1.627.03'
$ws.Range('E14').Value = '  -0.73%  '
Set-TextValue 'D15' '0.556'
$ws.Range('E15').Value = '  -0.05%  '
Set-TextValue 'D16' '0.0₃0760'
$ws.Range('E16').Value = '  -2.07%  '
Set-TextValue 'D17' '63.02'
$ws.Range('E17').Value = '  -0.17%  '
Set-TextValue 'D18' '25.761.87'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  -0.09%  '
Set-TextValue 'D20' '4.43'
$ws.Range('E20').Value = '  -0.24%  '
Set-TextValue 'D21' '191.81'
$ws.Range('E21').Value = '  -1.43%  '
Set-TextValue 'D22' '9.91'
$ws.Range('E22').Value = '  -0.20%  '
Set-TextValue 'D23' '6.26'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  +3.07%  '
Set-TextValue 'D26' '143.05'
$ws.Range('E26').Value = '  +2.07%  '
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  -1.00%  '
Set-TextValue 'D34' '1.56'
$ws.Range('E34').Value = '  -1.54%  '
Set-TextValue 'D35' '2.38'
$ws.Range('E35').Value = '  -0.36%  '
Set-TextValue 'D36' '0.902'
$ws.Range('E36').Value = '  +0.36%  '
Set-TextValue 'D37' '1.132.69'
$ws.Range('E37').Value = '  +1.82%  '
$ws.Range('E39').Value = '  -1.56%  '
Set-TextValue 'D40' '0.0154'
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('E42').Value = '  +0.90%  '
Set-TextValue 'D43' '100.37'
$ws.Range('E43').Value = '  +1.13%  '
Set-TextValue 'D44' '5.52'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('E45').Value = '  -0.33%  '
Set-TextValue 'D46' '1.765.05'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('E47').Value = '  +0.57%  '
Set-TextValue 'D48' '55.28'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D51' '1.42'
$ws.Range('E51').Value = '  +2.68%  '
